# Insert a new data row at row 15 (pushing the existing rows 15-58 down to
# 16-59), then populate the new row with its data. This mirrors the diff:
# dimension grows from A1:R58 to A1:R59, and every row from 15 downward is
# the previous row's content shifted down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(15).Insert()

$ws.Cells.Item(15, 1).Value = 4
$ws.Cells.Item(15, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(15, 3).Value = "Los Lagos"
$ws.Cells.Item(15, 4).Value = 44953
$ws.Cells.Item(15, 5).Value = 10
$ws.Cells.Item(15, 6).Value = 100112030
$ws.Cells.Item(15, 7).Value = "Poroto granado"
$ws.Cells.Item(15, 8).Value = "Sin especificar"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 80
$ws.Cells.Item(15, 11).Value = 45000
$ws.Cells.Item(15, 12).Value = 45000
$ws.Cells.Item(15, 13).Value = 45000
$ws.Cells.Item(15, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(15, 15).Value = "Región Metropolitana"
$ws.Cells.Item(15, 16).Value = 1800
$ws.Cells.Item(15, 17).Value = 25
$ws.Cells.Item(15, 18).Value = "Hortaliza"

$ws.Cells.Item(15, 4).NumberFormat = $ws.Cells.Item(16, 4).NumberFormat
